$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5674603174603174
$ws.Range("J2").Value = 0.03174603174603174
$ws.Range("P2").Value = 0.1507936507936508
$ws.Range("S2").Value = 0.05952380952380952
$ws.Range("B3").Value = 0.01408450704225352
$ws.Range("C3").Value = 0.007042253521126761
$ws.Range("J3").Value = 0.04225352112676056
$ws.Range("P3").Value = 0.7676056338028169
$ws.Range("S3").Value = 0.1690140845070423
$ws.Range("J4").Value = 0.03703703703703703
$ws.Range("P4").Value = 0.5555555555555556
$ws.Range("S4").Value = 0.4074074074074074
$ws.Range("B6").Value = 0.02816901408450704
$ws.Range("D6").Value = 0.0187793427230047
$ws.Range("F6").Value = 0.07042253521126761
$ws.Range("J6").Value = 0.2206572769953052
$ws.Range("O6").Value = 0.02816901408450704
$ws.Range("Q6").Value = 0.1596244131455399
$ws.Range("R6").Value = 0.07981220657276995
$ws.Range("S6").Value = 0.3943661971830986
$ws.Range("B7").Value = 0.078125
$ws.Range("D7").Value = 0.005208333333333333
$ws.Range("F7").Value = 0.08333333333333333
$ws.Range("J7").Value = 0.1510416666666667
$ws.Range("O7").Value = 0.015625
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.4166666666666667
$ws.Range("B8").Value = 0.06997742663656885
$ws.Range("D8").Value = 0.004514672686230248
$ws.Range("E8").Value = 0.004514672686230248
$ws.Range("F8").Value = 0.05417607223476298
$ws.Range("J8").Value = 0.0835214446952596
$ws.Range("O8").Value = 0.02483069977426636
$ws.Range("Q8").Value = 0.218961625282167
$ws.Range("R8").Value = 0.08803611738148984
$ws.Range("S8").Value = 0.4514672686230248
$ws.Range("B9").Value = 0.07262569832402235
$ws.Range("D9").Value = 0.00558659217877095
$ws.Range("F9").Value = 0.07262569832402235
$ws.Range("J9").Value = 0.106145251396648
$ws.Range("Q9").Value = 0.217877094972067
$ws.Range("R9").Value = 0.0782122905027933
$ws.Range("S9").Value = 0.446927374301676
$ws.Range("B10").Value = 0.10625
$ws.Range("D10").Value = 0.01640625
$ws.Range("E10").Value = 0.00078125
$ws.Range("F10").Value = 0.06328125
$ws.Range("J10").Value = 0.12109375
$ws.Range("O10").Value = 0.00703125
$ws.Range("Q10").Value = 0.24765625
$ws.Range("R10").Value = 0.0671875
$ws.Range("S10").Value = 0.3703125
$ws.Range("G11").Value = 0.1298245614035088
$ws.Range("J11").Value = 0.08421052631578947
$ws.Range("K11").Value = 0.1789473684210526
$ws.Range("L11").Value = 0.5929824561403508
$ws.Range("S11").Value = 0.01403508771929825
$ws.Range("G12").Value = 0.770949720670391
$ws.Range("J12").Value = 0.1452513966480447
$ws.Range("K12").Value = 0.0111731843575419
$ws.Range("L12").Value = 0.0446927374301676
$ws.Range("S12").Value = 0.02793296089385475
$ws.Range("G13").Value = 0.7352941176470589
$ws.Range("J13").Value = 0.2647058823529412
$ws.Range("F15").Value = 0.02475247524752475
$ws.Range("H15").Value = 0.1386138613861386
$ws.Range("I15").Value = 0.06930693069306931
$ws.Range("J15").Value = 0.4455445544554456
$ws.Range("K15").Value = 0.06435643564356436
$ws.Range("M15").Value = 0.004950495049504951
$ws.Range("O15").Value = 0.0297029702970297
$ws.Range("S15").Value = 0.2227722772277228
$ws.Range("F16").Value = 0.01257861635220126
$ws.Range("H16").Value = 0.1949685534591195
$ws.Range("I16").Value = 0.1006289308176101
$ws.Range("J16").Value = 0.4276729559748428
$ws.Range("K16").Value = 0.0880503144654088
$ws.Range("O16").Value = 0.06289308176100629
$ws.Range("S16").Value = 0.1132075471698113
$ws.Range("F17").Value = 0.01941747572815534
$ws.Range("H17").Value = 0.1844660194174757
$ws.Range("I17").Value = 0.08737864077669903
$ws.Range("J17").Value = 0.4077669902912621
$ws.Range("K17").Value = 0.1087378640776699
$ws.Range("M17").Value = 0.01359223300970874
$ws.Range("N17").Value = 0.001941747572815534
$ws.Range("O17").Value = 0.06407766990291262
$ws.Range("S17").Value = 0.112621359223301
$ws.Range("F18").Value = 0.02325581395348837
$ws.Range("H18").Value = 0.1279069767441861
$ws.Range("I18").Value = 0.06395348837209303
$ws.Range("J18").Value = 0.4941860465116279
$ws.Range("K18").Value = 0.1104651162790698
$ws.Range("M18").Value = 0.005813953488372093
$ws.Range("O18").Value = 0.0755813953488372
$ws.Range("S18").Value = 0.09883720930232558
$ws.Range("F19").Value = 0.01443464314354451
$ws.Range("H19").Value = 0.2197273456295108
$ws.Range("I19").Value = 0.07778668805132317
$ws.Range("J19").Value = 0.3809141940657578
$ws.Range("K19").Value = 0.1002405773857257
$ws.Range("M19").Value = 0.02085004009623095
$ws.Range("O19").Value = 0.06655974338412189
$ws.Range("S19").Value = 0.1194867682437851
